# PinMapping.xlsx edit: ETH -> ETH(SPI0), add SCK pin entry, reset selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ETH" usage label to "ETH(SPI0)" (column D, row 1).
$ws.Range("D1").Value = "ETH(SPI0)"

# Row 15 (pin 13) now carries the SCK entry in the Usage column (B).
$ws.Range("B15").Value = "SCK"

# Move the view back to the top and park the selection on E11.
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
